# bitaxeUltra BOM.xlsx - "updated BOM xlsx and added Gerbers"
#
# R11 (100k, DNP) was removed from the BOM entirely, so delete its row
# from the BOM sheet. Excel will automatically renumber/shift the rows
# below it and drop the now-unused shared strings for R11's Value/DK/PARTNO.

$wb  = $excel.ActiveWorkbook
$ws  = $wb.Worksheets.Item("BOM")

$ws.Rows.Item(29).Delete() | Out-Null

# Leave the selection where the author left it after deleting the row.
$ws.Range("I38").Select() | Out-Null

# Re-activate the sheet that was active/visible before the edit so the
# workbook keeps opening on the "DK Order" tab.
$ws2 = $wb.Worksheets.Item("DK Order")
$ws2.Activate() | Out-Null
